# Generate Report for Handback
#
# The CI run produced a new handback-verification result for
# 97880680-601d-4eab-b077-4a15054f3633 in both the zh-cn and de-de
# language sheets: the target xliff that was handed back turned out to
# be stale, so the status engine recorded the "not latest" error detail,
# a new handback datetime, the handback target-file name, and a link to
# the (outdated) handback markdown file in the "Latest Target File"
# column.

$wb = $excel.ActiveWorkbook

$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1b821f2e8b327ab98837eda4f31d878a9993fd/e2e/97880680-601d-4eab-b077-4a15054f3633.md"
$targetMdName = "97880680-601d-4eab-b077-4a15054f3633.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d52f6d26f4a9e27c9806a04019f7eabd38fc62cf/e2e/97880680-601d-4eab-b077-4a15054f3633.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1b821f2e8b327ab98837eda4f31d878a9993fd/e2e/97880680-601d-4eab-b077-4a15054f3633.md."

# ---- zh-cn sheet, row 7 (97880680-601d-4eab-b077-4a15054f3633) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("J7").Value = "97880680-601d-4eab-b077-4a15054f3633.999bd2878ff898b52b8377d45c194b9ac9616ebe.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-16 20:51:30"
$wsZhCn.Range("P7").Value = $errorDetail

# Latest Target File (I7) becomes a hyperlink to the handback markdown
# file, same as the other rows in this column; Excel applies the builtin
# Hyperlink style automatically.
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $targetMdUrl, $null, $null, $targetMdName)

# ---- de-de sheet, row 7 (97880680-601d-4eab-b077-4a15054f3633) ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G7").Value = "97880680-601d-4eab-b077-4a15054f3633.999bd2878ff898b52b8377d45c194b9ac9616ebe.de-de.xlf"
$wsDeDe.Range("J7").Value = "97880680-601d-4eab-b077-4a15054f3633.999bd2878ff898b52b8377d45c194b9ac9616ebe.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-16 20:51:38"
$wsDeDe.Range("P7").Value = $errorDetail

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $targetMdUrl, $null, $null, $targetMdName)
